$p = $ppt.ActivePresentation

# Slide 1 (index 1): textbox "The moon"
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(2).TextFrame.TextRange.Text = "The moon"

# Notes for slide 1: "chicken and dumplings"
$n1 = $s1.NotesPage
$n1.Shapes.Item(2).TextFrame.TextRange.Text = "chicken and dumplings"

# Slide 2 (index 2): textbox "Demonstration of simple table syntax, with alignment"
$s2 = $p.Slides.Item(2)
$s2.Shapes.Item(2).TextFrame.TextRange.Text = "Demonstration of simple table syntax, with alignment"

# Notes for slide 2: "foo bar"
$n2 = $s2.NotesPage
$n2.Shapes.Item(2).TextFrame.TextRange.Text = "foo bar"

# Notes for slide 3: merge runs within each paragraph, keep blank middle paragraph
$s3 = $p.Slides.Item(3)
$n3 = $s3.NotesPage
$n3.Shapes.Item(2).TextFrame.TextRange.Text = "Some notes inside a column`r`rSome notes outside the column"
